# fall 13 week 3 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.25

$ws.Range("B4").Value = 1.49
$ws.Range("D4").Value = 1.33
$ws.Range("F4").Value = 1.08
$ws.Range("G4").Value = 0.86

$ws.Range("D5").Value = 1.33
$ws.Range("F5").Value = 1.02
$ws.Range("G5").Value = 0.73

$ws.Range("D7").Value = 1.86
$ws.Range("E7").Value = 1.9
$ws.Range("F7").Value = 1.45
